$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current last row (row 46, "99"/"99999"),
# shifting it down to row 47, so we can add the new "93"/"99993" row
# in its place (new row 46). Excel's row Insert carries the formatting
# of the row above into the newly inserted row.
$ws.Rows("46:46").Insert()

# Populate the new row's values
$ws.Range("A46").Value = 93
$ws.Range("B46").Value = 99993
